$wb = $excel.ActiveWorkbook

# --- Update the "Count" sheet with the latest effort counts (row 7) ---
$count = $wb.Worksheets.Item("Count")
$count.Activate()

$count.Range("B7").Value = "Travis Thayer"
$count.Range("C7").Value = 0
$count.Range("D7").Value = 0
$count.Range("E7").Value = 1
$count.Range("F7").Value = 0
$count.Range("G7").Value = 1
$count.Range("H7").Value = 1
$count.Range("I7").Value = 1

# Match row 7's formatting to row 6's formatting
$count.Range("B6:I6").Copy()
$count.Range("B7:I7").PasteSpecial(-4122)  # xlPasteFormats

$count.Range("K7").Select()
